# Update PR log from #35
# Append a new row (row 12) to the PR log sheet with the data for PR #35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 35
$ws.Range("B12").Value = "time added"
$ws.Range("C12").Value = "riya-morankar"
$ws.Range("D12").Value = "N/A"
$ws.Range("E12").Value = "edit1 to main"

# Force the date column to be stored as plain text (matching the rest of
# the sheet, which stores dates as literal strings rather than date
# serials) instead of letting Excel auto-convert the string to a date.
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "2025-06-18"
